$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.335.22'
$ws.Range("E2").Value = '  +1.52%  '
$ws.Range("D3").Value = '3.399.19'
$ws.Range("E3").Value = '  +1.30%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = '''581.57'
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("D6").Value = '''179.04'
$ws.Range("E6").Value = '  +0.76%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("E8").Value = '  +0.64%  '
$ws.Range("E9").Value = '  +7.53%  '
$ws.Range("D10").Value = '''0.586'
$ws.Range("E10").Value = '  +0.59%  '
$ws.Range("D11").Value = '''48.50'
$ws.Range("E11").Value = '  +0.82%  '
$ws.Range("D12").Value = '''0.0000283'
$ws.Range("E12").Value = '  +3.27%  '
$ws.Range("D13").Value = '''680.83'
$ws.Range("E13").Value = '  -1.21%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '3.948.30'
$ws.Range("E14").Value = '  +1.34%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '''8.60'
$ws.Range("E15").Value = '  +1.86%  '
$ws.Range("D16").Value = '69.472.72'
$ws.Range("E16").Value = '  +1.73%  '
$ws.Range("D17").Value = '3.398.82'
$ws.Range("E17").Value = '  +1.23%  '
$ws.Range("E18").Value = '  +0.60%  '
$ws.Range("D19").Value = '''17.71'
$ws.Range("E19").Value = '  +1.26%  '
$ws.Range("D20").Value = '''11.28'
$ws.Range("E20").Value = '  +0.37%  '
$ws.Range("D21").Value = '''0.909'
$ws.Range("E21").Value = '  +1.53%  '
$ws.Range("E22").Value = '  -1.82%  '
$ws.Range("D23").Value = '''17.07'
$ws.Range("E23").Value = '  +0.32%  '
$ws.Range("D24").Value = '''101.22'
$ws.Range("E24").Value = '  +0.67%  '
$ws.Range("E25").Value = '  -0.42%  '
$ws.Range("D26").Value = '''2.70'
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("D27").Value = '''9.73'
$ws.Range("E27").Value = '  +2.28%  '
$ws.Range("D28").Value = '''33.53'
$ws.Range("E28").Value = '  +1.48%  '
$ws.Range("D29").Value = '''8.75'
$ws.Range("E29").Value = '  +2.43%  '
$ws.Range("D30").Value = '''6.89'
$ws.Range("E30").Value = '  -1.36%  '
$ws.Range("D31").Value = '''3.78'
$ws.Range("E31").Value = '  +12.10%  '
$ws.Range("D32").Value = '''555.77'
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("E33").Value = '  -0.46%  '
$ws.Range("E34").Value = '  +0.25%  '
$ws.Range("D35").Value = '''58.07'
$ws.Range("E36").Value = '  +0.17%  '
$ws.Range("D37").Value = '3.607.10'
$ws.Range("E37").Value = '  -2.85%  '
$ws.Range("D38").Value = '''0.140'
$ws.Range("E38").Value = '  +2.09%  '
$ws.Range("D39").Value = '''35.26'
$ws.Range("E39").Value = '  +1.26%  '
$ws.Range("D40").Value = '0.0₃0744'
$ws.Range("E40").Value = '  +10.19%  '
$ws.Range("D41").Value = '''3.30'
$ws.Range("E41").Value = '  +3.87%  '
$ws.Range("D42").Value = '''2.70'
$ws.Range("E42").Value = '  +3.28%  '
$ws.Range("D43").Value = '''0.0426'
$ws.Range("E43").Value = '  +3.76%  '
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("E45").Value = '  +1.36%  '
$ws.Range("E46").Value = '  +0.25%  '
$ws.Range("E47").Value = '  +3.95%  '
$ws.Range("E48").Value = '  -0.13%  '
$ws.Range("D49").Value = '''131.17'
$ws.Range("E49").Value = '  -0.55%  '
$ws.Range("D50").Value = '''2.64'
$ws.Range("E50").Value = '  +2.56%  '
$ws.Range("E51").Value = '  -0.53%  '
